# Adds the ranking data for 2025-08-05 as a brand-new worksheet at the end of
# the workbook, modeled on the existing "2025-08-04" sheet (same layout:
# header row "rank","title","author","latest_episode" followed by 50 ranked
# rows).

$wb = $excel.ActiveWorkbook

# Remember which sheet was active so the tab selection in the workbook is
# left exactly as it was found (adding/copying a sheet makes Excel jump to
# the new tab otherwise).
$activeSheetName = $wb.ActiveSheet.Name

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)

# Duplicate the most recent day's sheet so the new tab inherits its layout,
# column/page setup, and header formatting, then drop it at the very end of
# the tab strip (immediately after the sheet we copied from).
$lastSheet.Copy([System.Reflection.Missing]::Value, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "2025-08-05"

$data = @(
    ,('1','魔王の俺が奴隷エルフを嫁にしたんだが、どう愛でればいい？','原作／手島史詞 キャラクター原案／COMTA 漫画／板垣ハコ','第72話')
    ,('2','勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～','漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり','第５１話　英雄を倒す器用貧乏（１）')
    ,('3','異世界メイドの三ツ星グルメ ～現代ごはん作ったら王宮で大バズリしました～','モリタ Ｕ４ nima','第13話（１）　ドーナツの騎士様とお土産スイーツ（１）')
    ,('4','濁る瞳で何を願う ハイセルク戦記','トルトネン 創-taro 斎藤八呑','第32話 轍')
    ,('5','ふかふかダンジョン攻略記～俺の異世界転生冒険譚～','KAKERU','第67話「フカフカのだんじょん」（前半）')
    ,('6','回復術士のやり直し','月夜涙(原作) 羽賀ソウケン(漫画) しおこんぶ(キャラクター原案)','第72話-2')
    ,('7','ガチャを回して仲間を増やす 最強の美少女軍団を作り上げろ','漫画：晴野しゅー 原作：ちんくるり キャラクター原案：イセ川ヤスタカ','第72話前半')
    ,('8','大学入学時から噂されていた美少女三姉妹、生き別れていた義妹だった。','こすずめ(著者) 夏乃実(原作) ポメ(キャラクター原案)','第9話②')
    ,('9','いつでも自宅に帰れる俺は、異世界で行商人をはじめました','漫画／明地雫 原作／霜月緋色 キャラクター原案／いわさきたかし','第59話')
    ,('10','異世界のすみっこで快適ものづくり生活 ～女神さまのくれた工房はちょっとやりすぎ性能だった～','西山アラタ(漫画) 長田信織(原作) 東上文(キャラクター原案)','休載イラスト')
    ,('11','無能と呼ばれた『精霊たらし』～実は異能で、精霊界では伝説的ヒーローでした～＠COMIC','原作：佐藤謙羊 漫画：タバタグランドキャニオン','第27話「求職は総合人材ギルドで」③')
    ,('12','異世界料理道','原作／EDA キャラクター原案・漫画／こちも','第65話')
    ,('13','ポーション頼みで生き延びます！ 続','原作：FUNA 漫画：園心ふつう キャラクター原案：すきま','第65話　新たなる旅立ち Ⅰ')
    ,('14','器用貧乏、城を建てる～開拓学園の劣等生なのに、 上級職のスキルと魔法がすべて使えます～＠COMIC','原作：佐藤謙羊 漫画：スガン','第22話①「衝撃の資産ランキング」')
    ,('15','転生したらスライムだった件　クレイマンREVENGE','原作：伏瀬 漫画：カジカ航 キャラクター原案：みっつばー','第29話　侵入')
    ,('16','尾守つみきと奇日常。','森下みゆ','第60話 つみきさんとケイドロ。')
    ,('17','暴食のベルセルク～俺だけレベルという概念を突破する～','漫画：滝乃大祐 原作：一色一凛 キャラクター原案：fame','第73話前半')
    ,('18','創造錬金術師は自由を謳歌する 故郷を追放されたら、魔王のお膝元で超絶効果のマジックアイテム作り放題になりました','姫乃タカ(漫画) 千月さかき(原作) かぼちゃ(キャラクター原案)','第19話-3')
    ,('19','真の聖女である私は追放されました。だからこの国はもう終わりです','松もくば 鬱沢色素 ぷきゅのすけ','第５０話　真の聖女、再び国を守る（１）')
    ,('20','異世界転生の冒険者','しばの番茶 ケンイチ ネム','story.61〈前半〉')
    ,('21','ロメリア戦記～伯爵令嬢、魔王を倒した後も人類やばそうだから軍隊組織する～','漫画：上戸 亮 原作：有山リョウ(小学館「ガガガブックス」刊) キャラクター原案：コダマ','第14話「得手不得手」①')
    ,('22','クラス最安値で売られた俺は、実は最強パラメーター','カンブリア爆発太郎(漫画) RYOMA(原作) 黒井ススム(キャラクター原案)','第36話-2')
    ,('23','婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版','漫画/すたひろ 原作/Y.A','chapter67【35話②】')
    ,('24','え、社内システム全てワンオペしている私を解雇ですか？','漫画：伊於 原作：下城米雪 キャラクター原案：icchi','4巻発売告知漫画')
    ,('25','没落予定の貴族だけど、暇だったから魔法を極めてみた@COMIC','漫画：秋咲りお 原作：三木なずな キャラクター原案：かぼちゃ','第53話 ①')
    ,('26','アレクサンダー英雄戦記～最強の土魔術士～','マツオカヨシノリ るれくちぇ なんじゃもんじゃ','第11話（後編）')
    ,('27','神の庭付き楠木邸','安斎アキラ(著者) えんじゅ(原作) ox(キャラクター原案)','第33話')
    ,('28','姫様“拷問”の時間です','原作:春原ロビンソン　漫画:ひらけい','拷問146')
    ,('29','願ってもない追放後からのスローライフ？ 〜引退したはずが成り行きで美少女ギャルの師匠になったらなぜかめちゃくちゃ懐かれた〜','ヤミーゴ(漫画) シュガースプーン。（GA文庫/SBクリエイティブ）(原作) なたーしゃ(キャラクター原案)','第5話-1：魔法のお勉強')
    ,('30','いとこのこ','いぬちく(著者)','第36話')
    ,('31','ジゼルの錬金飴','漫画： katoson 原作：斯波 キャラクター原案：LINO','第9話')
    ,('32','迷宮ブラックカンパニー','安村洋平','第50話　落花流水（後編）')
    ,('33','Ｒｅ：ゼロから始める異世界生活 第四章 聖域と強欲の魔女','花鶏ハルノ(作画) 相川有(構成) 長月達平(原作) 大塚真一郎(キャラクター原案)','第62話①　エリオール大森林の永久凍土')
    ,('34','異世界たぬき','小林安曇','第3話①')
    ,('35','ポンコツスキルしか使えない悪役魔女だけど、テイムしたパリピなスライムたちと強く生きます！','漫画：鈴木イゾ 原作：雨傘ヒョウゴ キャラクター原案：朝日川日和','第9話')
    ,('36','魔物ノ森ノ少女ノヴァ','kamatama','第7話「山の民」前編①')
    ,('37','追放されたギルド職員は、世界最強の召喚士@COMIC','原作：月島秀一 漫画：あづち涼 キャラクター原案：チワワ丸','第10話①「下から来た者」')
    ,('38','スライムは最強たる可能性を秘めている～２回目の人生、ちゃんとスライムと向き合います～@COMIC','オサフネオウジ（漫画） 犬型大（原作） 風花風花（キャラクター原案）','第7話')
    ,('39','転生陰陽師・賀茂一樹～二度と地獄はご免なので、閻魔大王の神気で無双します〜@COMIC','芳井りょう（漫画） 赤野用介（原作） hakusai（キャラクター原案）','第12話 ①')
    ,('40','安達としまむら','柚原もけ(漫画) 入間人間(原作) のん(キャラクターデザイン)','第49話「最初の旅の端２」①')
    ,('41','不純な彼女達は懺悔しない','ポロロッカ(著者)','第29話')
    ,('42','ニャイト・オブ・ザ・リビングキャット','原作：ホークマン 作画：メカルーツ','Chapter16　トランスフォーミャー～キャット覚醒～（後編）')
    ,('43','アラフォーおっさんはスローライフの夢を見るか？','漫画：大関詠詞 原作：サイトウアユム キャラクター原案： ジョンディー','第15話')
    ,('44','となりの殺し屋ちゃん','ジェイ・加藤(著者)','最終話')
    ,('45','まったく最近の探偵ときたら','五十嵐正邦(著者)','第114話')
    ,('46','さんしょく弁当','兎月あい(著者)','第18話#4')
    ,('47','我輩は猫魔導師である ～キジトラ・ルークの快適チート猫生活～','原作：猫神信仰研究会 漫画：三國大和 キャラクター原案：ハム','第23話')
    ,('48','遊遊じてき。','カンケー','第2話')
    ,('49','千年英雄','原作/福島航平 作画/中村ゆきひろ','20話①')
    ,('50','異世界で最強のスキルを生み出せたので、ひたすら無双することにしました。　～俺だけがスキルの数値を勝手に操作～','漫画：星トマジロウ 原作：ヒゲ抜き地蔵 キャラクター原案：山椒魚','第10話 ①')
)

$rowCount = $data.Count
$colCount = 4

$values = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $values[$i,0] = [int]$data[$i][0]
    $values[$i,1] = [string]$data[$i][1]
    $values[$i,2] = [string]$data[$i][2]
    $values[$i,3] = [string]$data[$i][3]
}

# Header row (rank/title/author/latest_episode) - keep identical to the
# other day-sheets.
$newSheet.Range("A1").Value = "rank"
$newSheet.Range("B1").Value = "title"
$newSheet.Range("C1").Value = "author"
$newSheet.Range("D1").Value = "latest_episode"

# Write all 50 ranked rows in one shot.
$endRow = 1 + $rowCount
$dataRange = $newSheet.Range("A2:D" + $endRow)
$dataRange.Value = $values

# Restore whichever sheet was active before we added/copied the new tab.
$wb.Worksheets.Item($activeSheetName).Activate()
